$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"24.07000000000032"
$ws.Range("H2").Value = [double]"0.00220460212857343"
$ws.Range("I2").Value = [double]"0.00220460212857343"
$ws.Range("L2").Value = [double]"40.99236417316133"
$ws.Range("M2").Value = "[10.243250890823774, 71.74147745549888]"
$ws.Range("N2").Value = [double]"0.01011737500441634"
$ws.Range("O2").Value = [double]"0.01011737500441634"
$ws.Range("P2").Value = [double]"2.270500396288118"
$ws.Range("Q2").Value = "[1.6038160694001933, 2.937184723176043]"
$ws.Range("R2").Value = [double]"1.647697955853289e-08"
$ws.Range("S2").Value = [double]"1.647697955853289e-08"
$ws.Range("T2").Value = [double]"71.30515794146015"
$ws.Range("U2").Value = "[55.73277037085056, 86.87754551206974]"
$ws.Range("V2").Value = [double]"6.164180277323794e-12"
$ws.Range("W2").Value = [double]"6.164180277323794e-12"
$ws.Range("X2").Value = [double]"15.37203203203224"
$ws.Range("Y2").Value = [double]"12.81805805805823"
$ws.Range("Z2").Value = [double]"17.92600600600625"

# Row 3
$ws.Range("F3").Value = [double]"24.07000000000032"
$ws.Range("H3").Value = [double]"0.01994502239603024"
$ws.Range("I3").Value = [double]"0.01994502239603024"
$ws.Range("L3").Value = [double]"36.41525086760399"
$ws.Range("M3").Value = "[7.927196009408064, 64.90330572579991]"
$ws.Range("N3").Value = [double]"0.01339760020536662"
$ws.Range("O3").Value = [double]"0.01339760020536662"
$ws.Range("P3").Value = [double]"1.490605523324886"
$ws.Range("Q3").Value = "[0.3836579616996527, 2.5975530849501194]"
$ws.Range("R3").Value = [double]"0.009434009489280193"
$ws.Range("S3").Value = [double]"0.009434009489280193"
$ws.Range("T3").Value = [double]"59.09935375230709"
$ws.Range("U3").Value = "[41.34965441680203, 76.84905308781214]"
$ws.Range("V3").Value = [double]"2.782895225372783e-08"
$ws.Range("W3").Value = [double]"2.782895225372783e-08"
$ws.Range("X3").Value = [double]"18.35969969969995"
$ws.Range("Y3").Value = [double]"14.11913913913933"
$ws.Range("Z3").Value = [double]"22.60026026026057"

# Row 4
$ws.Range("B4").Value = [double]"1"
$ws.Range("F4").Value = [double]"24.07000000000032"
$ws.Range("H4").Value = [double]"1.456473276750803e-05"
$ws.Range("I4").Value = [double]"1.456473276750803e-05"
$ws.Range("L4").Value = [double]"57.34576836680959"
$ws.Range("M4").Value = "[32.206029797381674, 82.4855069362375]"
$ws.Range("N4").Value = [double]"3.503416093564304e-05"
$ws.Range("O4").Value = [double]"3.503416093564304e-05"
$ws.Range("P4").Value = [double]"1.62897396852804"
$ws.Range("Q4").Value = "[1.1258159859711165, 2.132131951084964]"
$ws.Range("R4").Value = [double]"5.25238688098284e-08"
$ws.Range("S4").Value = [double]"5.25238688098284e-08"
$ws.Range("T4").Value = [double]"75.02388961401846"
$ws.Range("U4").Value = "[59.79375375792777, 90.25402547010916]"
$ws.Range("V4").Value = [double]"6.652456363553938e-13"
$ws.Range("W4").Value = [double]"6.652456363553938e-13"
$ws.Range("X4").Value = [double]"17.82962962962987"
$ws.Range("Y4").Value = [double]"15.90210210210232"
$ws.Range("Z4").Value = [double]"19.75715715715742"

# Row 5
$ws.Range("F5").Value = [double]"24.07000000000032"
$ws.Range("H5").Value = [double]"0.03794772636350552"
$ws.Range("I5").Value = [double]"0.03794772636350552"
$ws.Range("L5").Value = [double]"32.57709067641478"
$ws.Range("M5").Value = "[4.111613457498613, 61.042567895330954]"
$ws.Range("N5").Value = [double]"0.02583243854576578"
$ws.Range("O5").Value = [double]"0.02583243854576578"
$ws.Range("P5").Value = [double]"1.515763422452733"
$ws.Range("Q5").Value = "[0.25786846606042335, 2.7736583788450426]"
$ws.Range("R5").Value = [double]"0.01929128690077686"
$ws.Range("S5").Value = [double]"0.01929128690077686"
$ws.Range("T5").Value = [double]"53.01244779703081"
$ws.Range("U5").Value = "[35.51484460114453, 70.51005099291709]"
$ws.Range("V5").Value = [double]"2.201986428751468e-07"
$ws.Range("W5").Value = [double]"2.201986428751468e-07"
$ws.Range("X5").Value = [double]"18.26332332332357"
$ws.Range("Y5").Value = [double]"13.44450450450469"
$ws.Range("Z5").Value = [double]"23.08214214214245"

# Row 6
$ws.Range("F6").Value = [double]"24.07000000000032"
$ws.Range("H6").Value = [double]"0.0003862778520917765"
$ws.Range("I6").Value = [double]"0.0003862778520917765"
$ws.Range("L6").Value = [double]"39.91789625512475"
$ws.Range("M6").Value = "[16.45720794047029, 63.3785845697792]"
$ws.Range("N6").Value = [double]"0.001314603705668382"
$ws.Range("O6").Value = [double]"0.001314603705668382"
$ws.Range("P6").Value = [double]"1.855395060678656"
$ws.Range("Q6").Value = "[1.1887107337907326, 2.5220793875665803]"
$ws.Range("R6").Value = [double]"1.198828945891961e-06"
$ws.Range("S6").Value = [double]"1.198828945891961e-06"
$ws.Range("T6").Value = [double]"55.93394826234631"
$ws.Range("U6").Value = "[42.80962166599308, 69.05827485869953]"
$ws.Range("V6").Value = [double]"4.94984053744929e-11"
$ws.Range("W6").Value = [double]"4.94984053744929e-11"
$ws.Range("X6").Value = [double]"16.96224224224247"
$ws.Range("Y6").Value = [double]"14.40826826826846"
$ws.Range("Z6").Value = [double]"19.51621621621648"

# Row 7
$ws.Range("F7").Value = [double]"24.07000000000032"
$ws.Range("H7").Value = [double]"0.001087398474985179"
$ws.Range("I7").Value = [double]"0.001087398474985179"
$ws.Range("L7").Value = [double]"40.65547843027822"
$ws.Range("M7").Value = "[13.347577436842982, 67.96337942371346]"
$ws.Range("N7").Value = [double]"0.004406848095940052"
$ws.Range("O7").Value = [double]"0.004406848095940052"
$ws.Range("P7").Value = [double]"2.06923720326535"
$ws.Range("Q7").Value = "[1.3899739268135027, 2.7485004797171966]"
$ws.Range("R7").Value = [double]"1.963970810425764e-07"
$ws.Range("S7").Value = [double]"1.963970810425764e-07"
$ws.Range("T7").Value = [double]"67.86799274876273"
$ws.Range("U7").Value = "[53.37172886144725, 82.36425663607821]"
$ws.Range("V7").Value = [double]"3.168132423070347e-12"
$ws.Range("W7").Value = [double]"3.168132423070347e-12"
$ws.Range("X7").Value = [double]"16.14304304304326"
$ws.Range("Y7").Value = [double]"13.54088088088106"
$ws.Range("Z7").Value = [double]"18.74520520520545"

# Row 8
$ws.Range("F8").Value = [double]"24.07000000000032"
$ws.Range("H8").Value = [double]"0.006375606939659662"
$ws.Range("I8").Value = [double]"0.006375606939659662"
$ws.Range("L8").Value = [double]"45.36153659770087"
$ws.Range("M8").Value = "[12.81895022599386, 77.90412296940788]"
$ws.Range("N8").Value = [double]"0.007356317704745097"
$ws.Range("O8").Value = [double]"0.007356317704745097"
$ws.Range("P8").Value = [double]"1.641552918091963"
$ws.Range("Q8").Value = "[0.6729738016698859, 2.6101320345140406]"
$ws.Range("R8").Value = [double]"0.001367073817477804"
$ws.Range("S8").Value = [double]"0.001367073817477804"
$ws.Range("T8").Value = [double]"72.25053356326242"
$ws.Range("U8").Value = "[53.046156892424165, 91.45491023410068]"
$ws.Range("V8").Value = [double]"1.430416896042175e-09"
$ws.Range("W8").Value = [double]"1.430416896042175e-09"
$ws.Range("X8").Value = [double]"17.78144144144168"
$ws.Range("Y8").Value = [double]"14.07095095095115"
$ws.Range("Z8").Value = [double]"21.49193193193222"

# Row 9
$ws.Range("F9").Value = [double]"24.07000000000032"
$ws.Range("H9").Value = [double]"6.588412777530639e-05"
$ws.Range("I9").Value = [double]"6.588412777530639e-05"
$ws.Range("L9").Value = [double]"42.86519469546109"
$ws.Range("M9").Value = "[21.46271903824851, 64.26767035267368]"
$ws.Range("N9").Value = [double]"0.0002098980143530049"
$ws.Range("O9").Value = [double]"0.0002098980143530049"
$ws.Range("P9").Value = [double]"1.779921363295117"
$ws.Range("Q9").Value = "[1.2138686329185768, 2.345974093671658]"
$ws.Range("R9").Value = [double]"9.981223292321317e-08"
$ws.Range("S9").Value = [double]"9.981223292321317e-08"
$ws.Range("T9").Value = [double]"53.52553020887785"
$ws.Range("U9").Value = "[41.054636392880816, 65.99642402487487]"
$ws.Range("V9").Value = [double]"4.051847746211479e-11"
$ws.Range("W9").Value = [double]"4.051847746211479e-11"
$ws.Range("X9").Value = [double]"17.2513713713716"
$ws.Range("Y9").Value = [double]"15.0829029029031"
$ws.Range("Z9").Value = [double]"19.41983983984011"

# Row 10
$ws.Range("F10").Value = [double]"23.84000000000029"
$ws.Range("H10").Value = [double]"3.041345521426209e-06"
$ws.Range("I10").Value = [double]"3.041345521426209e-06"
$ws.Range("L10").Value = [double]"64.27492713207273"
$ws.Range("M10").Value = "[37.81126440286411, 90.73858986128135]"
$ws.Range("N10").Value = [double]"1.317505622888682e-05"
$ws.Range("O10").Value = [double]"1.317505622888682e-05"
$ws.Range("P10").Value = [double]"1.754763464167271"
$ws.Range("Q10").Value = "[1.2516054816103477, 2.257921446724195]"
$ws.Range("R10").Value = [double]"9.38234112624059e-09"
$ws.Range("S10").Value = [double]"9.38234112624059e-09"
$ws.Range("T10").Value = [double]"57.54036624983351"
$ws.Range("U10").Value = "[41.92920183560663, 73.1515306640604]"
$ws.Range("V10").Value = [double]"2.408408805720796e-09"
$ws.Range("W10").Value = [double]"2.408408805720796e-09"
$ws.Range("X10").Value = [double]"17.18198198198219"
$ws.Range("Y10").Value = [double]"15.27287287287306"
$ws.Range("Z10").Value = [double]"19.09109109109132"

# Row 11
$ws.Range("B11").Value = [double]"0"
$ws.Range("F11").Value = [double]"23.84000000000029"
$ws.Range("H11").Value = [double]"0.1276146506652338"
$ws.Range("I11").Value = [double]"0.1276146506652338"
$ws.Range("L11").Value = [double]"27.07326374251548"
$ws.Range("M11").Value = "[-6.924854878263556, 61.07138236329452]"
$ws.Range("N11").Value = [double]"0.1157406523901334"
$ws.Range("O11").Value = [double]"0.1157406523901334"
$ws.Range("P11").Value = [double]"1.83023716155081"
$ws.Range("Q11").Value = "[0.018868424345884094, 3.6416058987557367]"
$ws.Range("R11").Value = [double]"0.04775827863544557"
$ws.Range("S11").Value = [double]"0.04775827863544557"
$ws.Range("T11").Value = [double]"61.18361085776085"
$ws.Range("U11").Value = "[42.68572461430883, 79.68149710121287]"
$ws.Range("V11").Value = [double]"3.238665180482769e-08"
$ws.Range("W11").Value = [double]"3.238665180482769e-08"
$ws.Range("X11").Value = [double]"16.89561561561582"
$ws.Range("Y11").Value = [double]"10.02282282282294"
$ws.Range("Z11").Value = [double]"23.7684084084087"

# Row 12
$ws.Range("F12").Value = [double]"23.84000000000029"
$ws.Range("H12").Value = [double]"0.0002729691434665282"
$ws.Range("I12").Value = [double]"0.0002729691434665282"
$ws.Range("L12").Value = [double]"53.02202709519492"
$ws.Range("M12").Value = "[22.678204610957238, 83.3658495794326]"
$ws.Range("N12").Value = [double]"0.001002534678232525"
$ws.Range("O12").Value = [double]"0.001002534678232525"
$ws.Range("P12").Value = [double]"1.968605606753965"
$ws.Range("Q12").Value = "[1.3396581285578097, 2.5975530849501194]"
$ws.Range("R12").Value = [double]"1.102705067257403e-07"
$ws.Range("S12").Value = [double]"1.102705067257403e-07"
$ws.Range("T12").Value = [double]"61.14225701945188"
$ws.Range("U12").Value = "[44.226176649746236, 78.05833738915751]"
$ws.Range("V12").Value = [double]"3.925687552808199e-09"
$ws.Range("W12").Value = [double]"3.925687552808199e-09"
$ws.Range("X12").Value = [double]"16.37061061061081"
$ws.Range("Y12").Value = [double]"13.98422422422439"
$ws.Range("Z12").Value = [double]"18.75699699699722"

# Row 13
$ws.Range("F13").Value = [double]"23.84000000000029"
$ws.Range("H13").Value = [double]"0.001047246681256331"
$ws.Range("I13").Value = [double]"0.001047246681256331"
$ws.Range("L13").Value = [double]"39.2895734616982"
$ws.Range("M13").Value = "[13.682780064504541, 64.89636685889185]"
$ws.Range("N13").Value = [double]"0.003424200727635274"
$ws.Range("O13").Value = [double]"0.003424200727635274"
$ws.Range("P13").Value = [double]"1.956026657190042"
$ws.Range("Q13").Value = "[1.2767633807381946, 2.6352899336418885]"
$ws.Range("R13").Value = [double]"6.18289395237781e-07"
$ws.Range("S13").Value = [double]"6.18289395237781e-07"
$ws.Range("T13").Value = [double]"60.67241740851041"
$ws.Range("U13").Value = "[46.73937896197848, 74.60545585504232]"
$ws.Range("V13").Value = [double]"2.679545474393308e-11"
$ws.Range("W13").Value = [double]"2.679545474393308e-11"
$ws.Range("X13").Value = [double]"16.41833833833854"
$ws.Range("Y13").Value = [double]"13.84104104104121"
$ws.Range("Z13").Value = [double]"18.99563563563586"

# Row 14
$ws.Range("F14").Value = [double]"23.84000000000029"
$ws.Range("H14").Value = [double]"0.00169466327391754"
$ws.Range("I14").Value = [double]"0.00169466327391754"
$ws.Range("L14").Value = [double]"41.07323112519221"
$ws.Range("M14").Value = "[13.072832398877424, 69.073629851507]"
$ws.Range("N14").Value = [double]"0.004968042088926206"
$ws.Range("O14").Value = [double]"0.004968042088926206"
$ws.Range("P14").Value = [double]"1.918289808498272"
$ws.Range("Q14").Value = "[1.1509738850989626, 2.6856057318975806]"
$ws.Range("R14").Value = [double]"8.179350519110784e-06"
$ws.Range("S14").Value = [double]"8.179350519110784e-06"
$ws.Range("T14").Value = [double]"55.14367604173018"
$ws.Range("U14").Value = "[39.95780250387979, 70.32954957958057]"
$ws.Range("V14").Value = [double]"3.498951350877633e-09"
$ws.Range("W14").Value = [double]"3.498951350877633e-09"
$ws.Range("X14").Value = [double]"16.56152152152172"
$ws.Range("Y14").Value = [double]"13.6501301301303"
$ws.Range("Z14").Value = [double]"19.47291291291315"
